$d = $word.ActiveDocument

# Locate the paragraph that currently holds the _GoBack bookmark up front
# (before we touch it) so we don't have to hard-code a paragraph index.
$bmOld = $d.Bookmarks("_GoBack")
$targetParaIndex = $bmOld.Range.Paragraphs.Item(1).Index

# 1) "D. Bridges and I. Harvey designed" -> "D. Bridges, I. Harvey and N. Qi designed"
$d.Content.Find.Execute("D. Bridges and I. Harvey designed", $true, $false, $false, $false, $false, $true, 1, $false, "D. Bridges, I. Harvey and N. Qi designed", 2)

# 2) "I. Harvey, N. Qi, E. Stephenson" -> "I. Harvey, E. Stephenson"
$d.Content.Find.Execute("I. Harvey, N. Qi, E. Stephenson", $true, $false, $false, $false, $false, $true, 1, $false, "I. Harvey, E. Stephenson", 2)

# 3) Append the new closing sentence, with a temporary trailing marker
#    character 'X' so we can precisely re-seat the _GoBack bookmark at the
#    true end of the paragraph afterwards (working around a quirk where a
#    zero-length bookmark placed directly at "end of paragraph text"
#    resolves incorrectly).
$d.Content.Find.Execute("edited and reviewed the manuscript.", $true, $false, $false, $false, $false, $true, 1, $false, "edited and reviewed the manuscript. All authors were involved in discussions. This manuscripts has been approved by all authors.X", 2)

# Relocate the _GoBack bookmark (currently mid-paragraph, left over from the
# original document) to the very end of the paragraph, matching where Word
# leaves it after this edit.
$bmOld = $d.Bookmarks("_GoBack")
$bmOld.Delete()

$targetPara = $d.Paragraphs($targetParaIndex)
$pr = $targetPara.Range
$lastPos = $pr.End - 1
$wrap = $d.Range($lastPos - 1, $lastPos)
$d.Bookmarks.Add("_GoBack", $wrap)
$bm = $d.Bookmarks("_GoBack")

# Delete the temporary 'X' marker via the bookmark's own range -- this
# leaves the bookmark collapsed exactly at the new end of the paragraph's
# text, with no stray character left behind.
$bmRange = $bm.Range
$bmRange.Text = ""
